$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: for numeric-looking "Price" column (D) values, force the cell
# to remain a text value (matching the source inlineStr type) by switching the
# NumberFormat to Text before assigning, then restoring the default "Normal" style
# so no stray style index is left on the cell.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.195.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.83%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.249.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.96%  "

# Row 7
$ws.Range("E7").Value = "  +2.13%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.69%  "

# Row 11
$ws.Range("E11").Value = "  +6.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0794"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.88%  "

# Row 13
$ws.Range("E13").Value = "  +2.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.94%  "

# Row 15
$ws.Range("E15").Value = "  +1.67%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.248.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.44%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.748"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.92%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.119.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.91%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.91%  "

# Row 22
$ws.Range("E22").Value = "  +1.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.10%  "

# Row 25
$ws.Range("E25").Value = "  +4.06%  "

# Row 26
$ws.Range("E26").Value = "  -0.19%  "

# Row 27
$ws.Range("E27").Value = "  +2.78%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.53%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.74%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0733"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.01%  "

# Row 36
$ws.Range("E36").Value = "  +7.00%  "

# Row 37
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.99%  "

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "

# Row 39
$ws.Range("E39").Value = "  +2.82%  "

# Row 40
$ws.Range("E40").Value = "  +5.15%  "

# Row 41
$ws.Range("E41").Value = "  +6.42%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.065.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0276"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.19%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.14%  "

# Row 48
$ws.Range("E48").Value = "  -3.25%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.468.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.94%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.73%  "
